$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Risks_Issues")

# Remove the "Supplier Lead Times" (R003) and "Currency Fluctuation" (R004) risk rows.
# These are currently rows 4 and 5; delete row 5 first so row indices stay valid.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Renumber the remaining risk IDs so they stay sequential (R001-R004).
$ws.Range("A4").Value = "R003"
$ws.Range("A5").Value = "R004"
